# Updates cryptocurrency market data table (Coin / Link / Price / Volume(1h))
# to reflect the latest refreshed values from the GitHub Actions data pull.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "70.785.62"
$ws.Range("E2").Value = "  +1.75%  "
$ws.Range("D3").Value = "3.629.20"
$ws.Range("E3").Value = "  +3.47%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.22%  "
$ws.Range("D5").Value = "'605.91"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "'199.93"
$ws.Range("E6").Value = "  +2.25%  "
$ws.Range("D7").Value = "'0.627"
$ws.Range("E7").Value = "  +0.94%  "
$ws.Range("D8").Value = "'1.00"
$ws.Range("E8").Value = "  +0.11%  "
$ws.Range("E9").Value = "  +9.35%  "
$ws.Range("D10").Value = "'0.647"
$ws.Range("E10").Value = "  +0.33%  "
$ws.Range("D11").Value = "'53.76"
$ws.Range("E11").Value = "  +1.13%  "
$ws.Range("E12").Value = "  +2.93%  "
$ws.Range("D13").Value = "'9.55"
$ws.Range("E13").Value = "  +1.18%  "
$ws.Range("D14").Value = "4.204.71"
$ws.Range("E14").Value = "  +3.63%  "
$ws.Range("D15").Value = "'680.24"
$ws.Range("E15").Value = "  +13.79%  "
$ws.Range("B16").Value = "WrappedBTC"
$ws.Range("C16").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D16").Value = "70.904.06"
$ws.Range("E16").Value = "  +1.69%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").Value = "'12.92"
$ws.Range("E17").Value = "  +2.18%  "
$ws.Range("D18").Value = "3.604.50"
$ws.Range("E18").Value = "  +2.77%  "
$ws.Range("D19").Value = "'19.00"
$ws.Range("E19").Value = "  +0.00%  "
$ws.Range("E20").Value = "  +0.03%  "
$ws.Range("D21").Value = "'0.999"
$ws.Range("E21").Value = "  +1.64%  "
$ws.Range("D22").Value = "'18.53"
$ws.Range("E22").Value = "  +2.98%  "
$ws.Range("D23").Value = "'5.39"
$ws.Range("E23").Value = "  +1.95%  "
$ws.Range("D24").Value = "'105.99"
$ws.Range("E24").Value = "  +4.00%  "
$ws.Range("E25").Value = "  +0.37%  "
$ws.Range("D26").Value = "'3.02"
$ws.Range("E26").Value = "  -3.28%  "
$ws.Range("D27").Value = "'10.50"
$ws.Range("E27").Value = "  -2.39%  "
$ws.Range("D28").Value = "'9.83"
$ws.Range("E28").Value = "  +3.67%  "
$ws.Range("D29").Value = "'34.29"
$ws.Range("E29").Value = "  +3.92%  "
$ws.Range("D30").Value = "'4.62"
$ws.Range("E30").Value = "  +7.72%  "
$ws.Range("D31").Value = "'7.19"
$ws.Range("E31").Value = "  +2.86%  "
$ws.Range("D32").Value = "'12.20"
$ws.Range("E32").Value = "  -0.88%  "
$ws.Range("D33").Value = "'0.115"
$ws.Range("E33").Value = "  +1.56%  "
$ws.Range("D34").Value = "'63.40"
$ws.Range("E34").Value = "  +0.38%  "
$ws.Range("D35").Value = "0.0₃0869"
$ws.Range("E35").Value = "  +7.33%  "
$ws.Range("D36").Value = "3.955.72"
$ws.Range("E36").Value = "  +5.78%  "
$ws.Range("E37").Value = "  +0.01%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "'514.01"
$ws.Range("E38").Value = "  +4.98%  "
$ws.Range("B39").Value = "Fetch.AI"
$ws.Range("C39").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D39").Value = "'3.01"
$ws.Range("E39").Value = "  -4.78%  "
$ws.Range("B40").Value = "InjectiveProtocol"
$ws.Range("C40").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D40").Value = "'36.60"
$ws.Range("E40").Value = "  +1.35%  "
$ws.Range("E41").Value = "  -0.13%  "
$ws.Range("E42").Value = "  -2.49%  "
$ws.Range("E43").Value = "  +2.11%  "
$ws.Range("B44").Value = "VeChain"
$ws.Range("C44").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D44").Value = "'0.0459"
$ws.Range("E44").Value = "  +2.31%  "
$ws.Range("B45").Value = "ThetaToken"
$ws.Range("C45").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D45").Value = "'3.07"
$ws.Range("E45").Value = "  +9.42%  "
$ws.Range("D46").Value = "'3.48"
$ws.Range("E46").Value = "  +5.98%  "
$ws.Range("E47").Value = "  +1.78%  "
$ws.Range("D48").Value = "'8.65"
$ws.Range("E48").Value = "  +3.13%  "
$ws.Range("E49").Value = "  -0.07%  "
$ws.Range("D50").Value = "'0.000247"
$ws.Range("E50").Value = "  +1.53%  "
$ws.Range("E51").Value = "  +2.05%  "
